$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("DW3").Value = 213.27
$ws.Range("DY3").Value = 913.27
$ws.Range("DW4").Value = 213.62
$ws.Range("DY4").Value = 1126.89
$ws.Range("DW5").Value = 211.15
$ws.Range("DY5").Value = 1338.04
$ws.Range("DW6").Value = 211.33
$ws.Range("DY6").Value = 1549.37
$ws.Range("DW7").Value = 212.91
$ws.Range("DY7").Value = 1762.28
$ws.Range("DW8").Value = 213.09
$ws.Range("DY8").Value = 1975.37
$ws.Range("DW9").Value = 211.33
$ws.Range("DY9").Value = 2186.7
$ws.Range("DW10").Value = 365.32
$ws.Range("DY10").Value = 2552.02
$ws.Range("DW11").Value = 211.5
$ws.Range("DY11").Value = 2763.52
$ws.Range("DW12").Value = -531.51
$ws.Range("DY12").Value = 2232.01
$ws.Range("DW13").Value = -532.92
$ws.Range("DY13").Value = 1699.09
$ws.Range("DW14").Value = 216.79
$ws.Range("DY14").Value = 1915.88
$ws.Range("DW15").Value = -528.15
$ws.Range("DY15").Value = 1387.73
$ws.Range("DW16").Value = 779.15
$ws.Range("DY16").Value = 2166.88
$ws.Range("DW17").Value = 776.6900000000001
$ws.Range("DY17").Value = 2943.57
$ws.Range("DW18").Value = 217.32
$ws.Range("DY18").Value = 3160.89
$ws.Range("DW19").Value = 353.5
$ws.Range("DY19").Value = 3514.39
$ws.Range("DW20").Value = -527.62
$ws.Range("DY20").Value = 2986.77
$ws.Range("DW21").Value = 918.16
$ws.Range("DY21").Value = 3904.93
$ws.Range("DW22").Value = 213.44
$ws.Range("DY22").Value = 4118.37
$ws.Range("DW23").Value = 354.56
$ws.Range("DY23").Value = 4472.93
$ws.Range("DW24").Value = -526.39
$ws.Range("DY24").Value = 3946.54
$ws.Range("DW25").Value = 214.15
$ws.Range("DY25").Value = 4160.690000000001
$ws.Range("DW26").Value = 212.74
$ws.Range("DY26").Value = 4373.43
$ws.Range("DW27").Value = 214.68
$ws.Range("DY27").Value = 4588.110000000001
$ws.Range("DW28").Value = 497.09
$ws.Range("DY28").Value = 5085.200000000001
$ws.Range("DY29").Value = 5085.200000000001
$ws.Range("DW30").Value = 214.68
$ws.Range("DY30").Value = 5299.880000000001
$ws.Range("DW31").Value = 211.68
$ws.Range("DY31").Value = 5511.560000000001
$ws.Range("DW32").Value = 212.74
$ws.Range("DY32").Value = 5724.300000000001
$ws.Range("DW33").Value = 500.27
$ws.Range("DY33").Value = 6224.570000000002
$ws.Range("DW34").Value = -523.91
$ws.Range("DY34").Value = 5700.660000000002
$ws.Range("DW35").Value = 493.74
$ws.Range("DY35").Value = 6194.400000000001
$ws.Range("DW36").Value = 212.74
$ws.Range("DY36").Value = 6407.140000000001
$ws.Range("DW37").Value = 360.91
$ws.Range("DY37").Value = 6768.050000000001
$ws.Range("DW38").Value = -527.27
$ws.Range("DY38").Value = 6240.780000000001
$ws.Range("DW39").Value = 492.68
$ws.Range("DY39").Value = 6733.460000000001
$ws.Range("DW40").Value = -530.27
$ws.Range("DY40").Value = 6203.190000000001
$ws.Range("DW41").Value = -531.34
$ws.Range("DY41").Value = 5671.85
$ws.Range("DW42").Value = 234.96
$ws.Range("DY42").Value = 5906.81
$ws.Range("DW43").Value = 211.15
$ws.Range("DY43").Value = 6117.96
$ws.Range("DW44").Value = -530.27
$ws.Range("DY44").Value = 5587.690000000001
$ws.Range("DW45").Value = 495.33
$ws.Range("DY45").Value = 6083.02
$ws.Range("DW46").Value = -528.86
$ws.Range("DY46").Value = 5554.160000000001
$ws.Range("DW47").Value = 200.57
$ws.Range("DY47").Value = 5754.73
$ws.Range("DW48").Value = 212.03
$ws.Range("DY48").Value = 5966.76
$ws.Range("DW49").Value = 212.56
$ws.Range("DY49").Value = 6179.320000000001
$ws.Range("DW50").Value = 213.8
$ws.Range("DY50").Value = 6393.120000000001
$ws.Range("DW51").Value = -544.76
$ws.Range("DY51").Value = 5848.360000000001
$ws.Range("DW52").Value = -497.06
$ws.Range("DY52").Value = 5351.3
$ws.Range("DW53").Value = 1063.51
$ws.Range("DY53").Value = 6414.81
$ws.Range("DW54").Value = 498.5
$ws.Range("DY54").Value = 6913.31
$ws.Range("DW55").Value = -528.15
$ws.Range("DY55").Value = 6385.160000000001
$ws.Range("DW56").Value = 493.57
$ws.Range("DY56").Value = 6878.73
$ws.Range("DW57").Value = 212.56
$ws.Range("DY57").Value = 7091.290000000001
$ws.Range("DW58").Value = -529.04
$ws.Range("DY58").Value = 6562.250000000001
$ws.Range("DW59").Value = 211.68
$ws.Range("DY59").Value = 6773.930000000001
